$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.250.38'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '3.469.61'
$ws.Range("E3").Value = '  +4.81%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '''243.75'
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").Value = '''646.87'
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '''1.42'
$ws.Range("E7").Value = '  +6.79%  '
$ws.Range("D8").Value = '''0.412'
$ws.Range("E8").Value = '  +1.53%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '''1.00'
$ws.Range("E10").Value = '  +4.75%  '
$ws.Range("D11").Value = '3.469.27'
$ws.Range("E11").Value = '  +4.85%  '
$ws.Range("D12").Value = '''43.42'
$ws.Range("E12").Value = '  +10.87%  '
$ws.Range("D13").Value = '''0.199'
$ws.Range("E13").Value = '  -1.53%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").Value = '''6.13'
$ws.Range("E14").Value = '  +3.55%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '95.877.01'
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("D16").Value = '4.117.56'
$ws.Range("E16").Value = '  +5.02%  '
$ws.Range("D17").Value = '''0.0000254'
$ws.Range("E17").Value = '  +3.22%  '
$ws.Range("D18").Value = '''8.61'
$ws.Range("E18").Value = '  +2.85%  '
$ws.Range("D19").Value = '3.469.14'
$ws.Range("E19").Value = '  +5.13%  '
$ws.Range("D20").Value = '''18.44'
$ws.Range("E20").Value = '  +11.68%  '
$ws.Range("D21").Value = '''11.98'
$ws.Range("E21").Value = '  +16.60%  '
$ws.Range("D22").Value = '''0.497'
$ws.Range("E22").Value = '  +8.03%  '
$ws.Range("D23").Value = '''516.12'
$ws.Range("E23").Value = '  +5.48%  '
$ws.Range("D24").Value = '''3.30'
$ws.Range("E24").Value = '  +0.31%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000194'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '''6.57'
$ws.Range("E26").Value = '  +3.92%  '
$ws.Range("D27").Value = '''92.29'
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").Value = '''12.46'
$ws.Range("E28").Value = '  +5.85%  '
$ws.Range("D29").Value = '3.647.55'
$ws.Range("E29").Value = '  +4.75%  '
$ws.Range("D30").Value = '''12.03'
$ws.Range("E30").Value = '  +14.05%  '
$ws.Range("D31").Value = '''1.00'
$ws.Range("E31").Value = '  +0.13%  '
$ws.Range("D32").Value = '''2.78'
$ws.Range("E32").Value = '  +15.40%  '
$ws.Range("D33").Value = '''0.140'
$ws.Range("E33").Value = '  +1.16%  '
$ws.Range("D34").Value = '''0.185'
$ws.Range("E34").Value = '  +1.89%  '
$ws.Range("B35").Value = 'EthereumClassic'
$ws.Range("C35").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D35").Value = '''31.18'
$ws.Range("E35").Value = '  +13.52%  '
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").Value = '''0.585'
$ws.Range("E36").Value = '  +10.15%  '
$ws.Range("E37").Value = '  +0.50%  '
$ws.Range("D38").Value = '''7.88'
$ws.Range("E38").Value = '  +6.81%  '
$ws.Range("D39").Value = '''1.48'
$ws.Range("E39").Value = '  +2.74%  '
$ws.Range("D40").Value = '''0.152'
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '''0.929'
$ws.Range("E42").Value = '  +15.07%  '
$ws.Range("D43").Value = '''513.62'
$ws.Range("E43").Value = '  +3.68%  '
$ws.Range("D44").Value = '''24.20'
$ws.Range("E44").Value = '  -1.03%  '
$ws.Range("D45").Value = '''1.73'
$ws.Range("E45").Value = '  +9.77%  '
$ws.Range("D46").Value = '''0.0425'
$ws.Range("E46").Value = '  +7.18%  '
$ws.Range("E47").Value = '  -0.33%  '
$ws.Range("D48").Value = '''5.60'
$ws.Range("E48").Value = '  +5.60%  '
$ws.Range("D49").Value = '''3.34'
$ws.Range("E49").Value = '  +7.89%  '
$ws.Range("D50").Value = '''2.19'
$ws.Range("E50").Value = '  +13.56%  '
$ws.Range("D51").Value = '''8.29'
$ws.Range("E51").Value = '  +1.44%  '
